$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 7) to the "AI Generated" sheet, mirroring the
# existing sample rows above it.
$ws.Range("A7").Value = "25 f3 address"
$ws.Range("B7").Value = " "
$ws.Range("C7").Value = "25 f3 firest"
$ws.Range("D7").Value = "25 f3 last"
$ws.Range("E7").Value = "{{ip_address}}"
$ws.Range("F7").Value = "25 city f3"

# G7 ("253") looks numeric, but the source data keeps it as text (like the
# other zip/state-ish columns in this sheet, e.g. G2="2502", G6="251").
# Force text formatting before assigning so Excel doesn't coerce it to a
# number, then drop back to the default "Normal" style so no unintended
# number-format style lingers on the cell.
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "253"
$ws.Range("G7").Style = "Normal"

$ws.Range("H7").Value = " "
